$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.207.22'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '1.682.02'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5278'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2688'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06367'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.48'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07629'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '1.690.33'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.520'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5761'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008239'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.44'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.51%  '
$ws.Range('D17').Value = '26.248.28'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.866'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('E20').Value = '  -0.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.241'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1262'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.725'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06396'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.377'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.021'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6126'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.415'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.745'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.176'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01636'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('D40').Value = '1.095.30'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8816'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '1.834.96'
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000110'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.107'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05265'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.015'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.02%  '
